$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4449
$ws.Range("L3").Value = 4735
$ws.Range("L4").Value = 1177
$ws.Range("L5").Value = 272
$ws.Range("L6").Value = 4076
$ws.Range("L7").Value = 14709

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 121
$ws.Range("L6").Value = 113
$ws.Range("L7").Value = 489
$ws.Range("L8").Value = 983
$ws.Range("L9").Value = 88
$ws.Range("L10").Value = 97
$ws.Range("L11").Value = 237
$ws.Range("L15").Value = 110
$ws.Range("L19").Value = 406
$ws.Range("L20").Value = 377
$ws.Range("L27").Value = 132
$ws.Range("L29").Value = 815
$ws.Range("L32").Value = 19
$ws.Range("L33").Value = 676
$ws.Range("L36").Value = 185
$ws.Range("L37").Value = 543
$ws.Range("L39").Value = 10
$ws.Range("L40").Value = 39
$ws.Range("L42").Value = 473
$ws.Range("L43").Value = 107
$ws.Range("L51").Value = 181
$ws.Range("L52").Value = 299
$ws.Range("L53").Value = 171
$ws.Range("L54").Value = 303
$ws.Range("L60").Value = 94
$ws.Range("L63").Value = 47
$ws.Range("L65").Value = 281
$ws.Range("L67").Value = 504
$ws.Range("L76").Value = 229
$ws.Range("L77").Value = 100
$ws.Range("L78").Value = 195
$ws.Range("L79").Value = 386
$ws.Range("L82").Value = 22
$ws.Range("L84").Value = 144
$ws.Range("L85").Value = 755
$ws.Range("L86").Value = 114
$ws.Range("L90").Value = 147
$ws.Range("L91").Value = 203
$ws.Range("L95").Value = 196
$ws.Range("L96").Value = 160
$ws.Range("L101").Value = 14709

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L3").Value = 44
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 166
$ws.Range("L3").Value = 158
$ws.Range("L7").Value = 489

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 72
$ws.Range("L7").Value = 237

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 227
$ws.Range("L3").Value = 307
$ws.Range("L6").Value = 157
$ws.Range("L7").Value = 755

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L6").Value = 80
$ws.Range("L7").Value = 299

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 171

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L3").Value = 325
$ws.Range("L6").Value = 266
$ws.Range("L7").Value = 983

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 187
$ws.Range("L6").Value = 208
$ws.Range("L7").Value = 676

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 159
$ws.Range("L7").Value = 543

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L5").Value = 5
$ws.Range("L7").Value = 281

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 148
$ws.Range("L6").Value = 115
$ws.Range("L7").Value = 504

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 73
$ws.Range("L6").Value = 147
$ws.Range("L7").Value = 303

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 242
$ws.Range("L3").Value = 306
$ws.Range("L7").Value = 815

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 143
$ws.Range("L6").Value = 118
$ws.Range("L7").Value = 406

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 107
$ws.Range("L7").Value = 229

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 158
$ws.Range("L7").Value = 473

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 73
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 386

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 114
$ws.Range("L3").Value = 123
$ws.Range("L7").Value = 377

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 185

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L2").Value = 3
$ws.Range("L6").Value = 10

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 41
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 62
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 57
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 39
